# ILV Arbeit mo nachmittag
# Moves the "Monday afternoon" entry (C7/D7 on the "Pflichtenheft techn."
# sheet) over to Tuesday (E7/F7), and fills in the new Wednesday/whatever-
# comes-next entry (E8/F8) with the old Grundstruktur note, bumping the
# "Auswertung" hours from 1.5h to 4.5h.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pflichtenheft techn.")

# Clear out the old Monday-afternoon entry (hours + description).
$ws.Range("C7").Value = $null
$ws.Range("D7").Value = $null

# New Tuesday entry: 8.5h, "Mindmap(2h),  Grobkonzept(5.5)"
$ws.Range("E7").Value = 8.5
$ws.Range("F7").Value = "Mindmap(2h),  Grobkonzept(5.5)"

# Row 8 (previously held the Mindmap note in F8) now gets the old
# Grundstruktur note, with the Auswertung time bumped from 1.5h to 4.5h,
# plus the corresponding 7.5h total.
$ws.Range("E8").Value = 7.5
$ws.Range("F8").Value = "-Grundstruktur (1h)`n- Kleinere Korrekturen(0.8h) + Batchdatei (0.2h)`n- KIS (1h)`n- Auswertung (4.5h)"

# The active-cell selection on this sheet moved from F8 to E7.
$ws.Range("E7").Select()
